# Track fuel-cost-driven changes in cargo dist transported down to the
# vehicle technology (engine type) level.
#
# - Rename existing "SoCDTtiNTY" sheet to "SoCDTtiNTY-psgr" (passenger).
# - Add a new "SoCDTtiNTY-frgt" sheet (freight) right after it, with the
#   same layout.
# - Both sheets get one header/value column per vehicle technology
#   (battery electric, natural gas, gasoline, diesel, plugin hybrid,
#   LPG, hydrogen) instead of the old single "passenger"/"freight"
#   column - every technology column simply repeats the old
#   passenger (resp. freight) value for now.

$wb = $excel.ActiveWorkbook

$vehicleTechs = @(
    "battery electric vehicle",
    "natural gas vehicle",
    "gasoline vehicle",
    "diesel vehicle",
    "plugin hybrid vehicle",
    "LPG vehicle",
    "hydrogen vehicle"
)

$rowLabels = @("LDVs", "HDVs", "aircraft", "rail", "ships", "motorbikes")

function Fill-SoCSheet($ws, $values) {
    $ws.Cells.Item(1, 1).Value = "Share that is New (dimensionless)"
    for ($j = 0; $j -lt $vehicleTechs.Count; $j++) {
        $ws.Cells.Item(1, $j + 2).Value = $vehicleTechs[$j]
    }

    for ($i = 0; $i -lt $rowLabels.Count; $i++) {
        $r = $i + 2
        $ws.Cells.Item($r, 1).Value = $rowLabels[$i]
        $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 8)).Value = $values[$i]
    }

    # Header row formatting.
    $ws.Rows.Item(1).RowHeight = 30
    $headerRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, 8))
    $headerRange.WrapText = $true
    $ws.Cells.Item(1, 1).Font.Bold = $true
    $techHeaderRange = $ws.Range($ws.Cells.Item(1, 2), $ws.Cells.Item(1, 8))
    $techHeaderRange.HorizontalAlignment = -4152

    # Column widths.
    $ws.Columns.Item(1).ColumnWidth = 18.333333333333332
    $colsBH = $ws.Range($ws.Cells.Item(1, 2), $ws.Cells.Item(1, 8)).EntireColumn
    $colsBH.ColumnWidth = 13.666666666666666
}

# --- Passenger sheet (rename in place) -------------------------------
$wsPsgr = $wb.Worksheets.Item("SoCDTtiNTY")
$wsPsgr.Name = "SoCDTtiNTY-psgr"
$psgrValues = @(
    0.076,
    0.0435,
    0.0416,
    0.029,
    0.029819999999999999,
    0.0587
)
Fill-SoCSheet $wsPsgr $psgrValues

# --- Freight sheet (new, inserted right after the passenger sheet) ---
$wsFrgt = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsPsgr)
$wsFrgt.Name = "SoCDTtiNTY-frgt"
$wsFrgt.Tab.Color = $wsPsgr.Tab.Color
$frgtValues = @(
    0.07,
    0.035,
    0.042,
    0.029,
    0.0303,
    0
)
Fill-SoCSheet $wsFrgt $frgtValues

# Restore the original active sheet/selection so the workbook opens on
# "About" like before.
$wb.Worksheets.Item("About").Activate()
